$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A19").Value = 'Are you always able to stop using drugs when you want to? (If never use drugs, answer “Yes”)'
$ws.Range("A41").Value = 'How often do you lose sleep due to late-night internet use?'
$ws.Range("A42").Value = 'How often do you find yourself anticipating when you will go online again?'
$ws.Range("A43").Value = 'How often do you snap, yell, or act annoyed if someone interrupts you while online?'
$ws.Range("A44").Value = 'How often do you find yourself saying ''just a few more minutes'' when online?'
$ws.Range("A45").Value = 'How often do you try to cut down the amount of time you spend online and fail?'
$ws.Range("A46").Value = 'How often do you choose to spend time online instead of spending time with others (family, friends)?'
$ws.Range("A47").Value = 'How often do you feel restless, moody, depressed, or irritable when attempting to cut down or stop internet use?'
$ws.Range("A48").Value = 'How often do you check your email before something else that you need to do?'
$ws.Range("A49").Value = 'How often do you fear that life without the Internet would be boring, empty, or joyless?'
$ws.Range("A50").Value = 'How often do you snap, yell, or act annoyed if you are unable to use the Internet?'
$ws.Range("A51").Value = 'How often do you fantasize about being online when you are offline?'
$ws.Range("A52").Value = 'How often do you lose track of time when online?'
$ws.Range("A53").Value = 'How often do your grades or school work suffer because of the amount of time you spend online?'
$ws.Range("A54").Value = 'How often do you find yourself staying online longer than you had planned?'
$ws.Range("A55").Value = 'How often do you use the Internet as a way of escaping from problems or relieving a dysphoric mood (e.g., feelings of helplessness, guilt, anxiety)?'
$ws.Range("A87").Value = 'I can depend on my friends for help if I need it.'
